# Apply updated statistics values (commit: "updated results and code")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellUpdates = @{
    "AA4" = 0.133
    "AB4" = 0.364
    "AC4" = 0.126
    "AE4" = 0.078
    "AF4" = 0.737
    "AG4" = 0.094
    "AH4" = 0.307
    "AI4" = 0.658
    "AJ4" = 0.172
    "AK4" = 0.415
    "AL4" = 0.703
    "AN4" = 0.339
    "AO4" = 0.699
    "B4" = 0.296
    "C4" = 0.051
    "E4" = 0.151
    "H4" = 0.189
    "J4" = 0.108
    "K4" = 0.342
    "L4" = 0.102
    "M4" = 0.32
    "N4" = 0.271
    "O4" = 0.02
    "P4" = 0.143
    "Q4" = 0.513
    "R4" = 0.217
    "S4" = 0.466
    "T4" = 0.283
    "W4" = 0.244
    "Y4" = 0.208
    "Z4" = 0.451
    "AA5" = 0.15
    "AB5" = 0.388
    "AC5" = 0.763
    "AD5" = 0.181
    "AE5" = 0.425
    "AF5" = 0.974
    "AH5" = 0.16
    "AI5" = 0.763
    "AJ5" = 0.181
    "AK5" = 0.425
    "AL5" = 0.921
    "AM5" = 0.073
    "AN5" = 0.27
    "AO5" = 0.886
    "B5" = 0.8159999999999999
    "C5" = 0.15
    "D5" = 0.388
    "E5" = 0.658
    "F5" = 0.225
    "G5" = 0.474
    "H5" = 0.8159999999999999
    "I5" = 0.15
    "J5" = 0.388
    "K5" = 0.658
    "L5" = 0.225
    "M5" = 0.474
    "N5" = 0.842
    "O5" = 0.133
    "P5" = 0.365
    "Q5" = 0.579
    "R5" = 0.244
    "S5" = 0.494
    "T5" = 0.579
    "U5" = 0.244
    "V5" = 0.494
    "W5" = 0.737
    "X5" = 0.194
    "Y5" = 0.44
    "Z5" = 0.8159999999999999
    "AC6" = 0.216
    "AF6" = 0.839
    "AI6" = 0.707
    "AL6" = 0.797
    "AO6" = 0.781
    "B6" = 0.434
    "E6" = 0.246
    "H6" = 0.307
    "K6" = 0.45
    "N6" = 0.41
    "Q6" = 0.544
    "T6" = 0.38
    "W6" = 0.367
    "Z6" = 0.581
    "AC7" = 0.379
    "AF7" = 0.915
    "AI7" = 0.739
    "AL7" = 0.867
    "AO7" = 0.84
    "B7" = 0.604
    "E7" = 0.394
    "H7" = 0.491
    "K7" = 0.555
    "N7" = 0.592
    "Q7" = 0.5639999999999999
    "T7" = 0.479
    "W7" = 0.525
    "Z7" = 0.702
    "AA8" = 0.15
    "AB8" = 0.387
    "AC8" = 0.655
    "AD8" = 0.177
    "AE8" = 0.42
    "AF8" = 0.893
    "AG8" = 0.046
    "AH8" = 0.215
    "AI8" = 0.753
    "AJ8" = 0.18
    "AK8" = 0.424
    "AL8" = 0.892
    "AM8" = 0.078
    "AN8" = 0.279
    "AO8" = 0.846
    "B8" = 0.751
    "C8" = 0.15
    "D8" = 0.387
    "E8" = 0.548
    "F8" = 0.191
    "G8" = 0.437
    "H8" = 0.703
    "I8" = 0.152
    "J8" = 0.39
    "K8" = 0.586
    "L8" = 0.204
    "M8" = 0.452
    "N8" = 0.749
    "O8" = 0.137
    "P8" = 0.371
    "Q8" = 0.55
    "R8" = 0.229
    "S8" = 0.479
    "T8" = 0.501
    "U8" = 0.206
    "V8" = 0.454
    "W8" = 0.662
    "X8" = 0.18
    "Y8" = 0.424
    "Z8" = 0.751
    "AA9" = 0.225
    "AB9" = 0.474
    "AC9" = 0.553
    "AD9" = 0.247
    "AE9" = 0.497
    "AF9" = 0.763
    "AG9" = 0.181
    "AH9" = 0.425
    "AI9" = 0.737
    "AJ9" = 0.194
    "AK9" = 0.44
    "AL9" = 0.842
    "AM9" = 0.133
    "AN9" = 0.365
    "AO9" = 0.781
    "B9" = 0.658
    "C9" = 0.225
    "D9" = 0.474
    "E9" = 0.421
    "F9" = 0.244
    "G9" = 0.494
    "H9" = 0.579
    "I9" = 0.244
    "J9" = 0.494
    "K9" = 0.5
    "N9" = 0.632
    "O9" = 0.233
    "P9" = 0.482
    "Q9" = 0.5
    "T9" = 0.395
    "U9" = 0.239
    "V9" = 0.489
    "W9" = 0.553
    "X9" = 0.247
    "Y9" = 0.497
    "Z9" = 0.658
    "AA10" = 0.15
    "AB10" = 0.388
    "AC10" = 0.658
    "AD10" = 0.225
    "AE10" = 0.474
    "AF10" = 0.974
    "AH10" = 0.16
    "AI10" = 0.763
    "AJ10" = 0.181
    "AK10" = 0.425
    "AL10" = 0.921
    "AM10" = 0.073
    "AN10" = 0.27
    "AO10" = 0.886
    "B10" = 0.8159999999999999
    "C10" = 0.15
    "D10" = 0.388
    "E10" = 0.579
    "F10" = 0.244
    "G10" = 0.494
    "H10" = 0.737
    "I10" = 0.194
    "J10" = 0.44
    "K10" = 0.658
    "L10" = 0.225
    "M10" = 0.474
    "N10" = 0.8159999999999999
    "O10" = 0.15
    "P10" = 0.388
    "Q10" = 0.579
    "R10" = 0.244
    "S10" = 0.494
    "T10" = 0.579
    "U10" = 0.244
    "V10" = 0.494
    "W10" = 0.737
    "X10" = 0.194
    "Y10" = 0.44
    "Z10" = 0.8159999999999999
    "AA11" = 0.15
    "AB11" = 0.388
    "AC11" = 0.711
    "AD11" = 0.206
    "AE11" = 0.454
    "AF11" = 0.974
    "AH11" = 0.16
    "AI11" = 0.763
    "AJ11" = 0.181
    "AK11" = 0.425
    "AL11" = 0.921
    "AM11" = 0.073
    "AN11" = 0.27
    "AO11" = 0.886
    "B11" = 0.8159999999999999
    "C11" = 0.15
    "D11" = 0.388
    "E11" = 0.658
    "F11" = 0.225
    "G11" = 0.474
    "H11" = 0.8159999999999999
    "I11" = 0.15
    "J11" = 0.388
    "K11" = 0.658
    "L11" = 0.225
    "M11" = 0.474
    "N11" = 0.842
    "O11" = 0.133
    "P11" = 0.365
    "Q11" = 0.579
    "R11" = 0.244
    "S11" = 0.494
    "T11" = 0.579
    "U11" = 0.244
    "V11" = 0.494
    "W11" = 0.737
    "X11" = 0.194
    "Y11" = 0.44
    "Z11" = 0.8159999999999999
    "AA12" = 0.32
    "AB12" = 0.5659999999999999
    "AC12" = 1.793
    "AD12" = 2.44
    "AE12" = 1.562
    "AF12" = 1.243
    "AG12" = 0.238
    "AH12" = 0.488
    "AI12" = 1.034
    "AJ12" = 0.033
    "AK12" = 0.182
    "AL12" = 1.086
    "AM12" = 0.078
    "AN12" = 0.28
    "AO12" = 1.121
    "B12" = 1.258
    "C12" = 0.32
    "D12" = 0.5659999999999999
    "E12" = 1.68
    "F12" = 1.098
    "G12" = 1.048
    "H12" = 1.613
    "I12" = 1.334
    "J12" = 1.155
    "K12" = 1.4
    "L12" = 0.5600000000000001
    "M12" = 0.748
    "N12" = 1.406
    "O12" = 0.616
    "P12" = 0.785
    "Z12" = 1.258
    "AA13" = 2.878
    "AB13" = 1.697
    "AC13" = 6.378
    "AD13" = 2.181
    "AE13" = 1.477
    "AF13" = 1.605
    "AG13" = 0.713
    "AH13" = 0.844
    "AI13" = 1.289
    "AJ13" = 0.364
    "AK13" = 0.603
    "AL13" = 1.579
    "AM13" = 0.717
    "AN13" = 0.847
    "AO13" = 1.491
    "B13" = 3.474
    "C13" = 1.46
    "D13" = 1.208
    "E13" = 4.594
    "F13" = 0.429
    "G13" = 0.655
    "H13" = 4.611
    "I13" = 0.627
    "J13" = 0.792
    "K13" = 2.265
    "L13" = 0.606
    "M13" = 0.779
    "N13" = 3.263
    "O13" = 0.72
    "P13" = 0.849
    "Z13" = 2.514
}

foreach ($cellRef in $cellUpdates.Keys) {
    $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
}
